$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = 45919
$ws.Range("B21").Value = "Sharekhan , Kanjur"
$ws.Range("C21").Value = "oops, hashmap internal working, linked hashmap internal working, design pattern - singletone 5 types `nproject - how you cut deployment time, why redis, `nprogram - list repeated words in string in given order"

$ws.Rows.Item(21).RowHeight = 60

$ws.Application.ActiveWindow.ScrollRow = 20
$null = $ws.Range("C26").Select()
